# "first login ask for genre and age"
# The "has_answered" flag (column D) is reset to FALSE for every user so
# that the app prompts everyone again on their next login to collect
# genre and age.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 keeps its formula form, but now evaluates FALSE() instead of TRUE().
$ws.Range("D2").Formula = "=FALSE()"

# D3:D13 become plain (non-formula) boolean FALSE values.
$ws.Range("D3:D13").Value = $false

# D8:D13 previously carried their own (slightly different) boolean number
# format; line them up with the format already used by D2:D7.
$ws.Range("D8:D13").NumberFormat = $ws.Range("D2").NumberFormat

# A8:C13 adopt the same direct formatting already used by A5:C7.
$ws.Range("A8:C13").HorizontalAlignment = $ws.Range("A5").HorizontalAlignment
$ws.Range("A8:C13").VerticalAlignment = $ws.Range("A5").VerticalAlignment
$ws.Range("A8:C13").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("A8:C13").Locked = $ws.Range("A5").Locked
$ws.Range("A8:C13").FormulaHidden = $ws.Range("A5").FormulaHidden

# Leave the selection on D4:D13 (the column being reset) and scroll the
# sheet back up so row 1 is visible again.
$ws.Range("D4:D13").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
